# Auto-generated Excel COM-interop edit script
# Updates specific cell values across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
# per scheduled runner data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3960
$ws.Range("I74").Value = 3950
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 3950
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -3014
$ws.Range("N74").Value = -5872

$ws.Range("H77").Value = 3960
$ws.Range("I77").Value = 3950
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 19750
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -15070
$ws.Range("N77").Value = -29360

$ws.Range("H107").Value = 750.75
$ws.Range("I107").Value = 495
$ws.Range("J107").Value = 1006.5
$ws.Range("K107").Value = 495
$ws.Range("L107").Value = 1006.5
$ws.Range("M107").Value = 1425
$ws.Range("N107").Value = -4846.5

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 722
$ws.Range("I25").Value = 722
$ws.Range("K25").Value = 722
$ws.Range("M25").Value = -320

$ws.Range("H102").Value = 1849.6666
$ws.Range("I102").Value = 1399.5
$ws.Range("K102").Value = 1399.5
$ws.Range("M102").Value = 222.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 406.75
$ws.Range("I5").Value = 257.14285
$ws.Range("K5").Value = 257.14285
$ws.Range("M5").Value = -144.14285

$ws.Range("H7").Value = 100499.5
$ws.Range("I7").Value = 200000
$ws.Range("K7").Value = 200000
$ws.Range("M7").Value = -199887

$ws.Range("H11").Value = 201.66667
$ws.Range("I11").Value = 202.5
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 202.5
$ws.Range("L11").Value = 200
$ws.Range("M11").Value = -62.5
$ws.Range("N11").Value = -480

$ws.Range("H105").Value = 2224
$ws.Range("I105").Value = 1688
$ws.Range("J105").Value = 2272.7273
$ws.Range("K105").Value = 1688
$ws.Range("L105").Value = 2272.7273
$ws.Range("M105").Value = 59
$ws.Range("N105").Value = -5766.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 10075975
$ws.Range("I6").Value = 13433800
$ws.Range("K6").Value = 13433800
$ws.Range("M6").Value = -13433687

$ws.Range("H7").Value = 28.583334
$ws.Range("I7").Value = 15.857142
$ws.Range("J7").Value = 46.4
$ws.Range("K7").Value = 15.857142
$ws.Range("L7").Value = 46.4
$ws.Range("M7").Value = 97.142858
$ws.Range("N7").Value = -272.4

$ws.Range("H17").Value = 7000
$ws.Range("I17").Value = 5000
$ws.Range("J17").Value = 8000
$ws.Range("K17").Value = 5000
$ws.Range("L17").Value = 8000
$ws.Range("M17").Value = -4826
$ws.Range("N17").Value = -8348

$ws.Range("H25").Value = 1699.6666
$ws.Range("I25").Value = 1699.6666
$ws.Range("K25").Value = 1699.6666
$ws.Range("M25").Value = -1525.6666

$ws.Range("H31").Value = 3648.4285
$ws.Range("I31").Value = 2423.75
$ws.Range("J31").Value = 5281.3335
$ws.Range("K31").Value = 2423.75
$ws.Range("L31").Value = 5281.3335
$ws.Range("M31").Value = -2128.75
$ws.Range("N31").Value = -5871.3335

$ws.Range("H34").Value = 3648.4285
$ws.Range("I34").Value = 2423.75
$ws.Range("J34").Value = 5281.3335
$ws.Range("K34").Value = 2423.75
$ws.Range("L34").Value = 5281.3335
$ws.Range("M34").Value = -2221.75
$ws.Range("N34").Value = -5685.3335

$ws.Range("H41").Value = 1479.8
$ws.Range("I41").Value = 1479.8
$ws.Range("K41").Value = 1479.8
$ws.Range("M41").Value = -1051.8

$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()

$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("M62").Value = -4376

$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 25000
$ws.Range("M65").Value = -21880

$ws.Range("H132").Value = 4544.5
$ws.Range("I132").Value = 1317
$ws.Range("K132").Value = 3951
$ws.Range("M132").Value = -1421

$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 2000
$ws.Range("I8").Value = 2000
$ws.Range("K8").Value = 6000
$ws.Range("M8").Value = -5861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 1036.7273
$ws.Range("I4").Value = 140
$ws.Range("J4").Value = 10004
$ws.Range("K4").Value = 140
$ws.Range("L4").Value = 10004
$ws.Range("M4").Value = -28
$ws.Range("N4").Value = -10228

$ws.Range("H113").Value = 1875
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 1875
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1875
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6215

$ws.Range("H134").Value = 61999.6
$ws.Range("J134").Value = 61999.6
$ws.Range("L134").Value = 185998.8
$ws.Range("N134").Value = -191068.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 1000
$ws.Range("K12").Value = 1000
$ws.Range("M12").Value = -830

$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

$ws.Range("H82").Value = 675
$ws.Range("J82").Value = 675
$ws.Range("L82").Value = 675
$ws.Range("N82").Value = -1397

$ws.Range("H85").Value = 675
$ws.Range("J85").Value = 675
$ws.Range("L85").Value = 675
$ws.Range("N85").Value = -3171

$ws.Range("H136").Value = 99004
$ws.Range("I136").Value = 99004
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 297012
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -294462
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 10178.75
$ws.Range("J4").Value = 10178.75
$ws.Range("L4").Value = 10178.75
$ws.Range("N4").Value = -10404.75

$ws.Range("H58").Value = 41137.2
$ws.Range("I58").Value = 33833
$ws.Range("J58").Value = 52093.5
$ws.Range("K58").Value = 33833
$ws.Range("L58").Value = 52093.5
$ws.Range("M58").Value = -33525
$ws.Range("N58").Value = -52709.5

$ws.Range("H136").Value = 1449.6666
$ws.Range("I136").Value = 1535.8182
$ws.Range("J136").Value = 1314.2858
$ws.Range("K136").Value = 1314.2858
$ws.Range("L136").Value = 3942.8574
$ws.Range("M136").Value = -2057.4546
$ws.Range("N136").Value = -9042.857400000001
